$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 5-7 (old FAPs-sending rows, now removed since ECs rows became FAPs rows)
$ws.Rows("5:7").Delete()

# Update row 2: FAPs | Wnt3 | Fzd1 | ECs
$ws.Cells.Item(2, 1).Value = "FAPs"
$ws.Cells.Item(2, 2).Value = "Wnt3"
$ws.Cells.Item(2, 3).Value = "Fzd1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.1741663333333333
$ws.Cells.Item(2, 8).Value = 0.5224989999999999
$ws.Cells.Item(2, 9).Value = 1
$ws.Cells.Item(2, 10).Value = 1
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.8596446666666667
$ws.Cells.Item(2, 14).Value = 2.578934
$ws.Cells.Item(2, 15).Value = 0.05286426382906832
$ws.Cells.Item(2, 16).Value = 0.05286426382906832
$ws.Cells.Item(2, 17).Value = 0.1497211595628889
$ws.Cells.Item(2, 18).Value = 1.347490436066
$ws.Cells.Item(2, 19).Value = 0.05286426382906832
$ws.Cells.Item(2, 20).Value = 0.05286426382906832

# Update row 3: FAPs | Wnt3 | Fzd1 | FAPs
$ws.Cells.Item(3, 1).Value = "FAPs"
$ws.Cells.Item(3, 2).Value = "Wnt3"
$ws.Cells.Item(3, 3).Value = "Fzd1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.1741663333333333
$ws.Cells.Item(3, 8).Value = 0.5224989999999999
$ws.Cells.Item(3, 9).Value = 1
$ws.Cells.Item(3, 10).Value = 1
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 10.435983
$ws.Cells.Item(3, 14).Value = 31.307949
$ws.Cells.Item(3, 15).Value = 0.6417658132713033
$ws.Cells.Item(3, 16).Value = 0.6417658132713032
$ws.Cells.Item(3, 17).Value = 1.817596893838999
$ws.Cells.Item(3, 18).Value = 16.358372044551
$ws.Cells.Item(3, 19).Value = 0.6417658132713033
$ws.Cells.Item(3, 20).Value = 0.6417658132713032

# Update row 4: FAPs | Wnt3 | Fzd1 | MuSCs
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Wnt3"
$ws.Cells.Item(4, 3).Value = "Fzd1"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.1741663333333333
$ws.Cells.Item(4, 8).Value = 0.5224989999999999
$ws.Cells.Item(4, 9).Value = 1
$ws.Cells.Item(4, 10).Value = 1
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 4.965729333333333
$ws.Cells.Item(4, 14).Value = 14.897188
$ws.Cells.Item(4, 15).Value = 0.3053699228996285
$ws.Cells.Item(4, 16).Value = 0.3053699228996284
$ws.Cells.Item(4, 17).Value = 0.8648628703124442
$ws.Cells.Item(4, 18).Value = 7.783765832811999
$ws.Cells.Item(4, 19).Value = 0.3053699228996285
$ws.Cells.Item(4, 20).Value = 0.3053699228996284
